$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B:E to remain plain text so Excel does not auto-convert
# numeric-looking strings (e.g. "240.45") into numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '35.345.53'
$ws.Range('E2').Value = '  -3.06%  '
$ws.Range('D3').Value = '1.975.63'
$ws.Range('E3').Value = '  -4.11%  '
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').Value = '240.45'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').Value = '0.631'
$ws.Range('E6').Value = '  -4.15%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '56.17'
$ws.Range('E8').Value = '  +7.47%  '
$ws.Range('D9').Value = '59.44'
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('E10').Value = '  -0.56%  '
$ws.Range('D11').Value = '0.0725'
$ws.Range('E11').Value = '  -2.65%  '
$ws.Range('D12').Value = '0.103'
$ws.Range('E12').Value = '  -5.60%  '
$ws.Range('D13').Value = '0.892'
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').Value = '14.17'
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').Value = '2.265.38'
$ws.Range('E15').Value = '  -4.70%  '
$ws.Range('E16').Value = '  -2.97%  '
$ws.Range('D17').Value = '1.985.51'
$ws.Range('E17').Value = '  -4.20%  '
$ws.Range('D18').Value = '17.11'
$ws.Range('E18').Value = '  +4.25%  '
$ws.Range('D19').Value = '35.219.49'
$ws.Range('E19').Value = '  -3.55%  '
$ws.Range('D20').Value = '69.83'
$ws.Range('E20').Value = '  -2.15%  '
$ws.Range('D21').Value = '0.0₃0833'
$ws.Range('E21').Value = '  -2.59%  '
$ws.Range('D22').Value = '231.63'
$ws.Range('E22').Value = '  -1.98%  '
$ws.Range('D23').Value = '5.02'
$ws.Range('E23').Value = '  -4.45%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').Value = '2.25'
$ws.Range('E25').Value = '  -5.15%  '
$ws.Range('D26').Value = '2.24'
$ws.Range('E26').Value = '  +5.77%  '
$ws.Range('D27').Value = '163.06'
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('D28').Value = '9.03'
$ws.Range('E28').Value = '  -4.04%  '
$ws.Range('D29').Value = '19.37'
$ws.Range('E29').Value = '  -4.45%  '
$ws.Range('E30').Value = '  -2.68%  '
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('E32').Value = '  -5.72%  '
$ws.Range('D33').Value = '0.0582'
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('D34').Value = '0.0891'
$ws.Range('E34').Value = '  +10.17%  '
$ws.Range('E35').Value = '  -6.95%  '
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').Value = '2.27'
$ws.Range('E37').Value = '  -3.34%  '
$ws.Range('E38').Value = '  -2.55%  '
$ws.Range('D39').Value = '4.83'
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D40').Value = '2.87'
$ws.Range('E40').Value = '  -1.40%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.18'
$ws.Range('E41').Value = '  -4.64%  '
$ws.Range('E42').Value = '  -4.00%  '
$ws.Range('E43').Value = '  -4.87%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = '0.0883'
$ws.Range('E44').Value = '  -5.78%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '90.20'
$ws.Range('E45').Value = '  -3.66%  '
$ws.Range('D46').Value = '1.357.32'
$ws.Range('E46').Value = '  -1.37%  '
$ws.Range('D47').Value = '7.39'
$ws.Range('E47').Value = '  -1.68%  '
$ws.Range('D48').Value = '15.35'
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('E50').Value = '  -3.72%  '
$ws.Range('D51').Value = '45.46'
$ws.Range('E51').Value = '  +2.16%  '
